$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# Fine-tune EE: mark Stuttgart EI ("Stuttgart EI" row) as chosen = "Yes"
$ws.Range("B4").Value = "Yes"

# Update the selected/active cell in the sheet view
$ws.Range("C3").Select()
